$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.120.94'
$ws.Range('E2').Value = '  -2.68%  '

# Row 3: 'Ethereum'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.711.33'
$ws.Range('E3').Value = '  -2.99%  '

# Row 4: 'TetherUSD'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.23%  '

# Row 5: 'BNB'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.14'
$ws.Range('E5').Value = '  -6.36%  '

# Row 6: 'USDC'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.10%  '

# Row 7: 'XRP'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4786'
$ws.Range('E7').Value = '  +7.16%  '

# Row 8: 'Cardano'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3441'
$ws.Range('E8').Value = '  -2.91%  '

# Row 9: 'OKB'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '41.92'
$ws.Range('E9').Value = '  -0.14%  '

# Row 10: 'Dogecoin'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07278'
$ws.Range('E10').Value = '  -1.75%  '

# Row 11: 'Polygon'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.048'
$ws.Range('E11').Value = '  -4.69%  '

# Row 12: 'BinanceUSD'
$ws.Range('E12').Value = '  -0.11%  '

# Row 13: 'Solana'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.86'
$ws.Range('E13').Value = '  -5.06%  '

# Row 14: 'Polkadot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.854'
$ws.Range('E14').Value = '  -2.80%  '

# Row 15: 'WrappedEther'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.705.65'
$ws.Range('E15').Value = '  -3.37%  '

# Row 16: 'Chainlink'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.841'
$ws.Range('E16').Value = '  -5.41%  '

# Row 17: 'Litecoin'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.15'
$ws.Range('E17').Value = '  -4.04%  '

# Row 18: 'ShibaInu'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001040'
$ws.Range('E18').Value = '  -1.98%  '

# Row 19: 'TRON'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06350'
$ws.Range('E19').Value = '  -1.18%  '

# Row 20: 'Dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9995'
$ws.Range('E20').Value = '  -0.14%  '

# Row 21: 'Avalanche'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.46'
$ws.Range('E21').Value = '  -3.75%  '

# Row 22: 'Uniswap'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.599'
$ws.Range('E22').Value = '  -2.81%  '

# Row 23: 'WrappedBTC'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.138.81'
$ws.Range('E23').Value = '  -2.77%  '

# Row 24: 'Cosmos'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.84'
$ws.Range('E24').Value = '  -3.59%  '

# Row 25: 'Toncoin'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.102'
$ws.Range('E25').Value = '  -0.25%  '

# Row 26: 'Monero'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.36'
$ws.Range('E26').Value = '  -3.47%  '

# Row 27: 'EthereumClassic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.58'
$ws.Range('E27').Value = '  -3.92%  '

# Row 28: 'WrappedliquidstakedEther2.0'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.874.88'
$ws.Range('E28').Value = '  -4.75%  '

# Row 29: 'LidoDAOToken'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.090'
$ws.Range('E29').Value = '  -2.79%  '

# Row 30: 'BitcoinCash'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.47'
$ws.Range('E30').Value = '  -4.01%  '

# Row 31: 'ImmutableX'
$ws.Range('E31').Value = '  -7.90%  '

# Row 32: 'Stellar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09198'
$ws.Range('E32').Value = '  +0.09%  '

# Row 33: 'HuobiToken'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.586'

# Row 34: 'Filecoin'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.312'
$ws.Range('E34').Value = '  -5.93%  '

# Row 35: 'VeChain'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02197'
$ws.Range('E35').Value = '  -3.68%  '

# Row 36: 'Hedera'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05828'
$ws.Range('E36').Value = '  -5.91%  '

# Row 37: 'Aptos'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.07'
$ws.Range('E37').Value = '  -6.52%  '

# Row 38: 'Algorand'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.1993'
$ws.Range('E38').Value = '  -5.12%  '

# Row 39: 'InternetComputer(DFINITY)'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.739'
$ws.Range('E39').Value = '  -4.23%  '

# Row 40: 'Frax'
$ws.Range('B40').Value = 'Frax'
$ws.Range('C40').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9993'
$ws.Range('E40').Value = '  -0.04%  '

# Row 41: 'WEMIXTOKEN'
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.400'
$ws.Range('E41').Value = '  +0.36%  '

# Row 42: 'TheSandbox'
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5865'
$ws.Range('E42').Value = '  -7.12%  '

# Row 43: 'TrustWalletToken'
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.107'
$ws.Range('E43').Value = '  -6.81%  '

# Row 44: 'FraxShare'
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.468'
$ws.Range('E44').Value = '  -5.06%  '

# Row 45: 'EnergySwap'
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.67'
$ws.Range('E45').Value = '  -4.20%  '

# Row 46: 'PancakeSwap'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.551'
$ws.Range('E46').Value = '  -5.34%  '

# Row 47: 'Decentraland'
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5626'
$ws.Range('E47').Value = '  -3.93%  '

# Row 48: 'Quant'
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '117.68'
$ws.Range('E48').Value = '  -3.86%  '

# Row 49: 'NEARProtocol'
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.840'
$ws.Range('E49').Value = '  -5.80%  '

# Row 50: 'Cronos'
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06629'
$ws.Range('E50').Value = '  -3.92%  '

# Row 51: 'EOS'
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.084'
$ws.Range('E51').Value = '  -4.42%  '
